$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: (blank) -> "false" (as literal text, not a Boolean).
# Assigning the bare word directly would auto-coerce to a Boolean TRUE/FALSE
# cell, so build it via a string-literal formula in a scratch cell and
# paste-special the computed value back in; that keeps the literal text
# "false" as a genuine text cell.
$ws.Range("Z1").Formula = "=""false"""
$ws.Range("Z1").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

# Date: updated publish timestamp
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# Case Sensitive: (blank) -> "true" (same text-not-Boolean trick as above)
$ws.Range("Z1").Formula = "=""true"""
$ws.Range("Z1").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
